$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New changelog entries (rows 27-33). Columns A (Date) and B (Version) contain
# values that look like a date / number ("2026-02-09", "1.0") so force those
# two columns to Text format first to stop Excel auto-converting them into a
# date serial / numeric value - the source data keeps them as literal text.
$ws.Range("A27:B33").NumberFormat = "@"

$rows = @(
    @{ Row = 27; Date = "2026-02-09"; Version = "1.0"; Category = "Feature";        Description = "Voice Capture: Added push-to-talk voice input for Day in the Life with real-time transcript preview (Web Speech API)" },
    @{ Row = 28; Date = "2026-02-09"; Version = "1.0"; Category = "Bug Fix";        Description = "Dashboard Mobile Scrolling: Fixed scrolling issues on mobile devices by removing height constraints" },
    @{ Row = 29; Date = "2026-02-09"; Version = "1.0"; Category = "Bug Fix";        Description = "Day in the Life Error Handling: Fixed streaming issues (infinite thinking, send button not activating, SSE error parsing)" },
    @{ Row = 30; Date = "2026-02-08"; Version = "1.0"; Category = "UI Enhancement"; Description = "AI Providers UI: Moved AI Providers management into Settings as tab with ListRow component for compact display" },
    @{ Row = 31; Date = "2026-02-08"; Version = "1.0"; Category = "Feature";        Description = "AI Providers Architecture: Implemented ai_providers table with CRUD operations, multi-provider support, real-time model fetching from Gemini/Anthropic/OpenAI APIs" },
    @{ Row = 32; Date = "2026-02-08"; Version = "1.0"; Category = "Feature";        Description = "Day in the Life: AI-powered chat interface for generating PACE-informed care observations. Supports Anthropic (Claude), OpenAI (GPT), Gemini with SSE streaming and custom prompts" },
    @{ Row = 33; Date = "2026-02-07"; Version = "1.0"; Category = "Feature";        Description = "Dashboard Enhancements: Enhanced company/personal dashboards with comprehensive data visibility, activity feeds, upcoming shifts, tabbed detail views for all features" }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Date
    $ws.Cells.Item($r, 2).Value = $entry.Version
    $ws.Cells.Item($r, 3).Value = $entry.Category
    $ws.Cells.Item($r, 4).Value = $entry.Description
    $ws.Cells.Item($r, 5).Value = "Claude"
}
